$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 29.75
$ws.Range("I42").Value = 10
$ws.Range("J42").Value = 36.333332
$ws.Range("K42").Value = 30
$ws.Range("L42").Value = 108.999996
$ws.Range("M42").Value = 200
$ws.Range("N42").Value = -568.999996

$ws.Range("H51").Value = 10571.071
$ws.Range("I51").Value = 9624.5
$ws.Range("J51").Value = 11833.167
$ws.Range("K51").Value = 9624.5
$ws.Range("L51").Value = 11833.167
$ws.Range("M51").Value = -9140.5
$ws.Range("N51").Value = -12801.167

$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()

$ws.Range("H129").Value = 3000
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 3000
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 9000
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -19000

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2422.0715
$ws.Range("J2").Value = 2150
$ws.Range("L2").Value = 2150
$ws.Range("N2").Value = -2376

$ws.Range("H23").Value = 1800
$ws.Range("J23").Value = 1800
$ws.Range("L23").Value = 1800
$ws.Range("N23").Value = -2318

$ws.Range("H74").Value = 834.5
$ws.Range("I74").Value = 823.5833
$ws.Range("K74").Value = 823.5833
$ws.Range("M74").Value = 50.41669999999999

$ws.Range("H77").Value = 834.5
$ws.Range("I77").Value = 823.5833
$ws.Range("K77").Value = 4117.9165
$ws.Range("M77").Value = 250.0834999999997

$ws.Range("H97").Value = 1018.0714
$ws.Range("I97").Value = 897.1667
$ws.Range("J97").Value = 1743.5
$ws.Range("K97").Value = 897.1667
$ws.Range("L97").Value = 1743.5
$ws.Range("M97").Value = -401.1667
$ws.Range("N97").Value = -2735.5

$ws.Range("H110").Value = 3944.2222
$ws.Range("I110").Value = 4187.25
$ws.Range("K110").Value = 4187.25
$ws.Range("M110").Value = -2142.25

$ws.Range("H116").Value = 2422.0715
$ws.Range("J116").Value = 2150
$ws.Range("L116").Value = 2150
$ws.Range("N116").Value = -6738

$ws.Range("H132").Value = 1663.091
$ws.Range("I132").Value = 1532.7778
$ws.Range("K132").Value = 4598.3334
$ws.Range("M132").Value = -2068.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2422.0715
$ws.Range("J3").Value = 2150
$ws.Range("L3").Value = 2150
$ws.Range("N3").Value = -2378

$ws.Range("H20").Value = 963.8889
$ws.Range("I20").Value = 897
$ws.Range("J20").Value = 1198
$ws.Range("K20").Value = 897
$ws.Range("L20").Value = 1198
$ws.Range("M20").Value = -650
$ws.Range("N20").Value = -1692

$ws.Range("H94").Value = 864.7895
$ws.Range("I94").Value = 595.5333000000001
$ws.Range("J94").Value = 1874.5
$ws.Range("K94").Value = 595.5333000000001
$ws.Range("L94").Value = 1874.5
$ws.Range("M94").Value = -144.5333000000001
$ws.Range("N94").Value = -2776.5

$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()

$ws.Range("H107").Value = 1693
$ws.Range("I107").Value = 897.4
$ws.Range("K107").Value = 897.4
$ws.Range("M107").Value = 1022.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2863
$ws.Range("I31").Value = 2720.125
$ws.Range("J31").Value = 3148.75
$ws.Range("K31").Value = 2720.125
$ws.Range("L31").Value = 3148.75
$ws.Range("M31").Value = -2425.125
$ws.Range("N31").Value = -3738.75

$ws.Range("H34").Value = 2863
$ws.Range("I34").Value = 2720.125
$ws.Range("J34").Value = 3148.75
$ws.Range("K34").Value = 2720.125
$ws.Range("L34").Value = 3148.75
$ws.Range("M34").Value = -2518.125
$ws.Range("N34").Value = -3552.75

$ws.Range("H74").Value = 34996
$ws.Range("J74").Value = 34996
$ws.Range("L74").Value = 34996
$ws.Range("N74").Value = -36744

$ws.Range("H77").Value = 34996
$ws.Range("J77").Value = 34996
$ws.Range("L77").Value = 104988
$ws.Range("N77").Value = -113724

$ws.Range("H99").Value = 3712.7
$ws.Range("I99").Value = 3175
$ws.Range("J99").Value = 4071.1667
$ws.Range("K99").Value = 3175
$ws.Range("L99").Value = 4071.1667
$ws.Range("M99").Value = -1677
$ws.Range("N99").Value = -7067.1667

$ws.Range("H126").Value = 3712.7
$ws.Range("I126").Value = 3175
$ws.Range("J126").Value = 4071.1667
$ws.Range("K126").Value = 9525
$ws.Range("L126").Value = 12213.5001
$ws.Range("M126").Value = -7055
$ws.Range("N126").Value = -17153.5001

$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 59.6
$ws.Range("I2").Value = 26.666666
$ws.Range("J2").Value = 109
$ws.Range("K2").Value = 159.999996
$ws.Range("L2").Value = 654
$ws.Range("M2").Value = -46.99999600000001
$ws.Range("N2").Value = -880

$ws.Range("H5").Value = 389.7143
$ws.Range("I5").Value = 413
$ws.Range("J5").Value = 250
$ws.Range("K5").Value = 1239
$ws.Range("L5").Value = 750
$ws.Range("M5").Value = -1127
$ws.Range("N5").Value = -974

$ws.Range("H11").Value = 11424.272
$ws.Range("J11").Value = 2227.4285
$ws.Range("L11").Value = 6682.2855
$ws.Range("N11").Value = -6962.2855

$ws.Range("H37").Value = 150000
$ws.Range("J37").Value = 150000
$ws.Range("L37").Value = 450000
$ws.Range("N37").Value = -450224

$ws.Range("H132").Value = 868.75
$ws.Range("I132").Value = 868.75
$ws.Range("K132").Value = 7818.75
$ws.Range("M132").Value = -5288.75

$ws.Range("H135").Value = 389.7143
$ws.Range("I135").Value = 413
$ws.Range("J135").Value = 250
$ws.Range("K135").Value = 3717
$ws.Range("L135").Value = 2250
$ws.Range("M135").Value = -1182
$ws.Range("N135").Value = -7320

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 1261
$ws.Range("J39").Value = 1261
$ws.Range("L39").Value = 1261
$ws.Range("N39").Value = -2325

$ws.Range("H70").Value = 337.5
$ws.Range("I70").Value = 337.5
$ws.Range("K70").Value = 337.5
$ws.Range("M70").Value = -67.5

$ws.Range("H73").Value = 337.5
$ws.Range("I73").Value = 337.5
$ws.Range("K73").Value = 337.5
$ws.Range("M73").Value = 598.5

$ws.Range("H97").Value = 2612.85
$ws.Range("I97").Value = 2544.7058
$ws.Range("J97").Value = 2999
$ws.Range("K97").Value = 2544.7058
$ws.Range("L97").Value = 2999
$ws.Range("M97").Value = -2048.7058
$ws.Range("N97").Value = -3991

$ws.Range("H102").Value = 1196
$ws.Range("I102").Value = 1136.5
$ws.Range("K102").Value = 1136.5
$ws.Range("M102").Value = 485.5

$ws.Range("H113").Value = 2554
$ws.Range("J113").Value = 3000
$ws.Range("L113").Value = 3000
$ws.Range("N113").Value = -7340

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()

$ws.Range("H132").Value = 1732.5
$ws.Range("J132").Value = 1700
$ws.Range("L132").Value = 5100
$ws.Range("N132").Value = -10160

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2247.5
$ws.Range("I40").Value = 2247.5
$ws.Range("K40").Value = 2247.5
$ws.Range("M40").Value = -2111.5

$ws.Range("H68").Value = 2268.3333
$ws.Range("I68").Value = 534
$ws.Range("J68").Value = 4002.6667
$ws.Range("K68").Value = 534
$ws.Range("L68").Value = 4002.6667
$ws.Range("M68").Value = 215
$ws.Range("N68").Value = -5500.6667

$ws.Range("H71").Value = 2268.3333
$ws.Range("I71").Value = 534
$ws.Range("J71").Value = 4002.6667
$ws.Range("K71").Value = 2670
$ws.Range("L71").Value = 20013.3335
$ws.Range("M71").Value = 1074
$ws.Range("N71").Value = -27501.3335

$ws.Range("H82").Value = 1575.375
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 1575.375
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 1575.375
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -2297.375

$ws.Range("H85").Value = 1575.375
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 1575.375
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 1575.375
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -4071.375

$ws.Range("H93").Value = 2233
$ws.Range("I93").Value = 2233
$ws.Range("K93").Value = 2233
$ws.Range("M93").Value = -985

$ws.Range("H132").Value = 5542.5
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()
